$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.435.77"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.97"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.021"
$ws.Range("E4").Value = "  +1.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.12"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.018"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5118"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3962"
$ws.Range("E8").Value = "  +2.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08460"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.93"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.268"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.13"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.48"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.235"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.021"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.93"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06777"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.018"
$ws.Range("E21").Value = "  +1.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.949"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.440.59"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.068.67"
$ws.Range("E26").Value = "  -1.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.86"
$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  +1.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.377"
$ws.Range("E29").Value = "  -3.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.51"
$ws.Range("E30").Value = "  +2.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1053"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.043"
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.788"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02435"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06465"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2184"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.872"
$ws.Range("E38").Value = "  -6.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.257"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.184"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6382"
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.997"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.26"
$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6019"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.96"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.707"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.990"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.208"
$ws.Range("E48").Value = "  -5.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.206"
$ws.Range("E49").Value = "  -1.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.76"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06868"
